# This workbook's data rows (2-26) got re-shuffled: the values in columns
# D (Fecha), H (Variedad), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), N (Unidad de comercializacion), O (Origen),
# P (Precio $/Kg) and Q (Kg o Unidades) move between rows according to a
# fixed permutation, while the rest of the row (Mercado ID, Mercado, Region,
# Codreg, Categoria ID, Categoria, Calidad, Clasificacion) stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get shuffled between rows: D=4 H=8 J=10 K=11 L=12 M=13 N=14 O=15 P=16 Q=17
$cols = @(4, 8, 10, 11, 12, 13, 14, 15, 16, 17)

# Snapshot the current (pre-edit) values for every data row before writing
# anything back, so the permutation can be applied safely in one pass.
$original = @{}
for ($r = 2; $r -le 26; $r++) {
    $rowvals = @{}
    foreach ($c in $cols) {
        $rowvals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $original[$r] = $rowvals
}

# Target row -> source row (i.e. row $target ends up with the values that
# row $source originally had).
$mapping = @{}
$mapping[2] = 23
$mapping[3] = 25
$mapping[4] = 12
$mapping[5] = 7
$mapping[6] = 21
$mapping[7] = 6
$mapping[8] = 16
$mapping[9] = 5
$mapping[10] = 13
$mapping[11] = 14
$mapping[12] = 22
$mapping[13] = 24
$mapping[14] = 4
$mapping[15] = 26
$mapping[16] = 20
$mapping[17] = 19
$mapping[18] = 2
$mapping[19] = 3
$mapping[20] = 9
$mapping[21] = 8
$mapping[22] = 17
$mapping[23] = 18
$mapping[24] = 10
$mapping[25] = 11
$mapping[26] = 15

for ($target = 2; $target -le 26; $target++) {
    $source = $mapping[$target]
    $src = $original[$source]
    foreach ($c in $cols) {
        $ws.Cells.Item($target, $c).Value2 = $src[$c]
    }
}
